$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new test case result for Test Case No.13 in column G (통과 여부 / Pass status)
$ws.Range("G14").Value = "추후 코드 수정 필요"
$ws.Range("G14").HorizontalAlignment = -4108
$ws.Range("G14").VerticalAlignment = -4108
$ws.Range("G14").WrapText = $true

# Update the view state: selection moved to G15, top-left cell to A10
$ws.Range("G15").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
